$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for "2022-Q3" and push the
#    existing "2022-Q1"/"2021-Q4" summary rows down one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Excel carries the header row's formatting into the freshly inserted row;
# strip it back to the plain (unstyled) look the data rows actually use.
$summary.Range("B2:D2").ClearFormats()

# Clone the style (bold/centered/bordered) used by the index column from the
# row right below (which just got pushed down and still carries it).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0

# Renumber the index column on the two rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q3" detail sheet, right after "总计". Duplicate the
#    "2022-Q1" sheet so formatting/styles come along for free, then replace
#    its data with the single 2022-Q3 holding and drop the two extra rows.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Rows.Item(3).Delete()
$q3.Rows.Item(3).Delete()

$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "002952"
$q3.Range("C2").Value = "建信多因子量化股票"

$q3.Range("D2:G2").NumberFormat = "@"
$q3.Range("D2").Value = "0.09"
$q3.Range("E2").Value = "91.26"
$q3.Range("F2").Value = "3.35"
$q3.Range("G2").Value = "0.0030"

$q3.Range("H2").Value = 6

Write-Output "done"
